{"js": "// The edit removes three paragraphs that immediately follow the\n// \"LOB1004: C\u00e1lculo II (Requisito fraco)\" paragraph:\n//   1. an empty paragraph\n//   2. \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   3. \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github\n//      pages. Original theme under Creative Commons Attribution\"\n// The paragraph that follows those three (another empty paragraph, then the\n// page-break paragraph) is left untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst anchorText = \"LOB1004: C\u00e1lculo II (Requisito fraco)\";\nconst targetTexts = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\nconst items = paragraphs.items;\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === anchorText) {\n    anchorIndex = i;\n    break;\n  }\n}\n\n// Only delete when the full 3-paragraph sequence (blank, \"Ver no\n// Jupiter...\", \"\u00a9 2020...\") is found right after the anchor - this avoids\n// accidentally deleting an unrelated blank paragraph if the script were to\n// run again on an already-edited document.\nif (\n  anchorIndex !== -1 &&\n  items[anchorIndex + 1] && items[anchorIndex + 1].text === \"\" &&\n  items[anchorIndex + 2] && items[anchorIndex + 2].text === targetTexts[0] &&\n  items[anchorIndex + 3] && items[anchorIndex + 3].text === targetTexts[1]\n) {\n  const toDelete = [\n    items[anchorIndex + 1],\n    items[anchorIndex + 2],\n    items[anchorIndex + 3]\n  ];\n\n  for (const p of toDelete) {\n    p.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# The edit removes three paragraphs that immediately follow the\n# \"LOB1004: C\u00e1lculo II (Requisito fraco)\" paragraph:\n#   1. an empty paragraph\n#   2. \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n#   3. \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github\n#      pages. Original theme under Creative Commons Attribution\"\n# The paragraph that follows those three (another empty paragraph, then the\n# page-break paragraph) is left untouched.\n\n$d = $word.ActiveDocument\n\n$anchorText = \"LOB1004: C\u00e1lculo II (Requisito fraco)\"\n$targetText1 = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$targetText2 = \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n\n$anchorIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.TrimEnd(\"`r\", \"`a\")\n    if ($t -eq $anchorText) {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -ge 1) {\n    $deleteAt = $anchorIndex + 1\n\n    $p1 = $d.Paragraphs.Item($deleteAt)\n    $t1 = $p1.Range.Text.TrimEnd(\"`r\", \"`a\")\n    if ($t1 -eq \"\") {\n        $p1.Range.Delete()\n    }\n\n    $p2 = $d.Paragraphs.Item($deleteAt)\n    $t2 = $p2.Range.Text.TrimEnd(\"`r\", \"`a\")\n    if ($t2 -eq $targetText1) {\n        $p2.Range.Delete()\n    }\n\n    $p3 = $d.Paragraphs.Item($deleteAt)\n    $t3 = $p3.Range.Text.TrimEnd(\"`r\", \"`a\")\n    if ($t3 -eq $targetText2) {\n        $p3.Range.Delete()\n    }\n}\n"}
